$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") contains a date value (serial 45179 = 2023-09-10)
# for every data row (rows 2-246). Update it to 45180 (2023-09-11) for all rows.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 246) { $lastRow = 246 }

$ws.Range("C2:C$lastRow").Value = 45180
